# Auto update Excel log
# Appends newly-logged sensor readings to the PIR, Humidity and Temperature
# sheets of the SeniorConnect master log.

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param($Sheet, $StartRow, $Rows)

    $endRow = $StartRow + $Rows.Count - 1
    $target = $Sheet.Range($Sheet.Cells.Item($StartRow, 1), $Sheet.Cells.Item($endRow, 6))
    # Force plain text so date-/time-/percent-looking strings are stored
    # verbatim instead of being auto-converted to numeric date/time/percent
    # values by Excel.
    $target.NumberFormat = "@"

    $r = $StartRow
    foreach ($row in $Rows) {
        $Sheet.Cells.Item($r, 1).Value = $row[0]
        $Sheet.Cells.Item($r, 2).Value = $row[1]
        $Sheet.Cells.Item($r, 3).Value = $row[2]
        $Sheet.Cells.Item($r, 4).Value = $row[3]
        $Sheet.Cells.Item($r, 5).Value = $row[4]
        $Sheet.Cells.Item($r, 6).Value = $row[5]
        $r = $r + 1
    }
}

# ---------------------------------------------------------------------------
# PIR sheet: append rows 99-111
# ---------------------------------------------------------------------------
$pir = $wb.Worksheets.Item("PIR")
$pirRows = @(
    ,@("2026-02-06","09:48:15","09:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-06","09:48:16","09:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-06","09:48:20","09:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-06","09:48:25","09:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-06","09:48:30","09:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-06","09:48:35","09:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-06","09:48:40","09:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-06","09:48:46","09:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-06","09:48:50","09:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-06","09:48:56","09:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-06","09:49:01","09:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-06","09:49:06","09:00","Bathroom","No Motion","Inactive")
    ,@("2026-02-06","09:49:11","09:00","Bathroom","No Motion","Inactive")
)
Add-LogRows $pir 99 $pirRows

# ---------------------------------------------------------------------------
# Humidity sheet: append rows 36-47
# ---------------------------------------------------------------------------
$humidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    ,@("2026-02-06","09:48:15","09:00","Bathroom","72.5%","Active")
    ,@("2026-02-06","09:48:16","09:00","Bathroom","72.5%","Active")
    ,@("2026-02-06","09:48:20","09:00","Bathroom","71.6%","Active")
    ,@("2026-02-06","09:48:25","09:00","Bathroom","72.5%","Active")
    ,@("2026-02-06","09:48:30","09:00","Bathroom","71.4%","Active")
    ,@("2026-02-06","09:48:35","09:00","Bathroom","72.4%","Active")
    ,@("2026-02-06","09:48:45","09:00","Bathroom","72.6%","Active")
    ,@("2026-02-06","09:48:50","09:00","Bathroom","71.5%","Active")
    ,@("2026-02-06","09:48:55","09:00","Bathroom","72.4%","Active")
    ,@("2026-02-06","09:49:00","09:00","Bathroom","71.4%","Active")
    ,@("2026-02-06","09:49:05","09:00","Bathroom","72.2%","Active")
    ,@("2026-02-06","09:49:10","09:00","Bathroom","71.2%","Active")
)
Add-LogRows $humidity 36 $humidityRows

# ---------------------------------------------------------------------------
# Temperature sheet: append rows 36-47
# ---------------------------------------------------------------------------
$temperature = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    ,@("2026-02-06","09:48:15","09:00","Bathroom","27.5C","Active")
    ,@("2026-02-06","09:48:16","09:00","Bathroom","27.5C","Active")
    ,@("2026-02-06","09:48:20","09:00","Bathroom","27.5C","Active")
    ,@("2026-02-06","09:48:25","09:00","Bathroom","27.5C","Active")
    ,@("2026-02-06","09:48:30","09:00","Bathroom","27.5C","Active")
    ,@("2026-02-06","09:48:35","09:00","Bathroom","27.5C","Active")
    ,@("2026-02-06","09:48:45","09:00","Bathroom","27.5C","Active")
    ,@("2026-02-06","09:48:50","09:00","Bathroom","27.5C","Active")
    ,@("2026-02-06","09:48:55","09:00","Bathroom","27.5C","Active")
    ,@("2026-02-06","09:49:00","09:00","Bathroom","27.5C","Active")
    ,@("2026-02-06","09:49:05","09:00","Bathroom","27.5C","Active")
    ,@("2026-02-06","09:49:11","09:00","Bathroom","27.5C","Active")
)
Add-LogRows $temperature 36 $temperatureRows
